$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4627.3887
$ws.Range("I74").Value = 4429.364
$ws.Range("J74").Value = 4938.5713
$ws.Range("K74").Value = 4429.364
$ws.Range("L74").Value = 4938.5713
$ws.Range("M74").Value = -3493.364
$ws.Range("N74").Value = -6810.5713
$ws.Range("H77").Value = 4627.3887
$ws.Range("I77").Value = 4429.364
$ws.Range("J77").Value = 4938.5713
$ws.Range("K77").Value = 22146.82
$ws.Range("L77").Value = 24692.8565
$ws.Range("M77").Value = -17466.82
$ws.Range("N77").Value = -34052.85649999999
$ws.Range("H92").Value = 384.6207
$ws.Range("I92").Value = 375.15384
$ws.Range("K92").Value = 375.15384
$ws.Range("M92").Value = 872.8461600000001
$ws.Range("H98").Value = 1915.6177
$ws.Range("I98").Value = 1133.2903
$ws.Range("J98").Value = 9999.666999999999
$ws.Range("K98").Value = 1133.2903
$ws.Range("L98").Value = 9999.666999999999
$ws.Range("M98").Value = 364.7097000000001
$ws.Range("N98").Value = -12995.667
$ws.Range("H103").Value = 922.5
$ws.Range("J103").Value = 1200
$ws.Range("L103").Value = 3600
$ws.Range("N103").Value = -4772
$ws.Range("H122").Value = 1915.6177
$ws.Range("I122").Value = 1133.2903
$ws.Range("J122").Value = 9999.666999999999
$ws.Range("K122").Value = 3399.8709
$ws.Range("L122").Value = 29999.001
$ws.Range("M122").Value = -949.8708999999999
$ws.Range("N122").Value = -34899.001
$ws.Range("H138").Value = 3825.3777
$ws.Range("I138").Value = 1913.7646
$ws.Range("J138").Value = 4986
$ws.Range("K138").Value = 5741.293799999999
$ws.Range("L138").Value = 14958
$ws.Range("M138").Value = -601.2937999999995
$ws.Range("N138").Value = -25238
$ws.Range("H141").Value = 1712253.9
$ws.Range("I141").Value = 2973.75
$ws.Range("K141").Value = 8921.25
$ws.Range("M141").Value = -3741.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 17860220
$ws.Range("I2").Value = 22728552
$ws.Range("J2").Value = 9666.666999999999
$ws.Range("K2").Value = 22728552
$ws.Range("L2").Value = 9666.666999999999
$ws.Range("M2").Value = -22728439
$ws.Range("N2").Value = -9892.666999999999
$ws.Range("H43").Value = 6000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 6000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 6000
$ws.Range("M43").Value = ""
$ws.Range("N43").Value = -6626
$ws.Range("H45").Value = 1580.2195
$ws.Range("I45").Value = 1058
$ws.Range("J45").Value = 3004.4546
$ws.Range("K45").Value = 1058
$ws.Range("L45").Value = 3004.4546
$ws.Range("M45").Value = -681
$ws.Range("N45").Value = -3758.4546
$ws.Range("H61").Value = 3037.4783
$ws.Range("I61").Value = 1628
$ws.Range("K61").Value = 1628
$ws.Range("M61").Value = -1416
$ws.Range("H97").Value = 511.72
$ws.Range("I97").Value = 469.69565
$ws.Range("K97").Value = 469.69565
$ws.Range("M97").Value = 26.30435
$ws.Range("H116").Value = 17860220
$ws.Range("I116").Value = 22728552
$ws.Range("J116").Value = 9666.666999999999
$ws.Range("K116").Value = 22728552
$ws.Range("L116").Value = 9666.666999999999
$ws.Range("M116").Value = -22726258
$ws.Range("N116").Value = -14254.667
$ws.Range("H122").Value = 4028.9167
$ws.Range("I122").Value = 2724.5
$ws.Range("J122").Value = 5333.3335
$ws.Range("K122").Value = 8173.5
$ws.Range("L122").Value = 16000.0005
$ws.Range("M122").Value = -5723.5
$ws.Range("N122").Value = -20900.0005
$ws.Range("H136").Value = 3037.4783
$ws.Range("I136").Value = 1628
$ws.Range("K136").Value = 4884
$ws.Range("M136").Value = -2334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 17860220
$ws.Range("I3").Value = 22728552
$ws.Range("J3").Value = 9666.666999999999
$ws.Range("K3").Value = 22728552
$ws.Range("L3").Value = 9666.666999999999
$ws.Range("M3").Value = -22728438
$ws.Range("N3").Value = -9894.666999999999
$ws.Range("H94").Value = 721.8
$ws.Range("I94").Value = 756
$ws.Range("J94").Value = 585
$ws.Range("K94").Value = 756
$ws.Range("L94").Value = 585
$ws.Range("M94").Value = -305
$ws.Range("N94").Value = -1487

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 29687
$ws.Range("J57").Value = 29687
$ws.Range("L57").Value = 29687
$ws.Range("N57").Value = -30807

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1158.6538
$ws.Range("I5").Value = 406.8421
$ws.Range("K5").Value = 1220.5263
$ws.Range("M5").Value = -1108.5263
$ws.Range("H113").Value = 633.4048
$ws.Range("I113").Value = 490.9655
$ws.Range("J113").Value = 951.1539
$ws.Range("K113").Value = 1472.8965
$ws.Range("L113").Value = 2853.4617
$ws.Range("M113").Value = 697.1034999999999
$ws.Range("N113").Value = -7193.4617
$ws.Range("H131").Value = 916.9
$ws.Range("J131").Value = 1112.96
$ws.Range("L131").Value = 3338.88
$ws.Range("N131").Value = -13418.88
$ws.Range("H135").Value = 1158.6538
$ws.Range("I135").Value = 406.8421
$ws.Range("K135").Value = 3661.5789
$ws.Range("M135").Value = -1126.5789

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2462.5
$ws.Range("I97").Value = 1451.6666
$ws.Range("J97").Value = 5495
$ws.Range("K97").Value = 1451.6666
$ws.Range("L97").Value = 5495
$ws.Range("M97").Value = -955.6666
$ws.Range("N97").Value = -6487
$ws.Range("H113").Value = 1911
$ws.Range("I113").Value = 1001.6667
$ws.Range("J113").Value = 3275
$ws.Range("K113").Value = 1001.6667
$ws.Range("L113").Value = 3275
$ws.Range("M113").Value = 1168.3333
$ws.Range("N113").Value = -7615
$ws.Range("H126").Value = 2792.3076
$ws.Range("I126").Value = 1518.75
$ws.Range("J126").Value = 4830
$ws.Range("K126").Value = 4556.25
$ws.Range("L126").Value = 14490
$ws.Range("M126").Value = -2086.25
$ws.Range("N126").Value = -19430
$ws.Range("H132").Value = 2588.8635
$ws.Range("I132").Value = 1964.9032
$ws.Range("K132").Value = 5894.7096
$ws.Range("M132").Value = -3364.7096

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2393.4375
$ws.Range("I7").Value = 1800
$ws.Range("J7").Value = 2855
$ws.Range("K7").Value = 1800
$ws.Range("L7").Value = 2855
$ws.Range("M7").Value = -1688
$ws.Range("N7").Value = -3079
$ws.Range("H40").Value = 4316.3335
$ws.Range("I40").Value = 1300
$ws.Range("J40").Value = 4919.6
$ws.Range("K40").Value = 1300
$ws.Range("L40").Value = 4919.6
$ws.Range("M40").Value = -1164
$ws.Range("N40").Value = -5191.6
$ws.Range("H61").Value = 100002360
$ws.Range("I61").Value = 111112620
$ws.Range("K61").Value = 111112620
$ws.Range("M61").Value = -111112418
$ws.Range("H68").Value = 1898.1818
$ws.Range("I68").Value = 1088
$ws.Range("J68").Value = 10000
$ws.Range("K68").Value = 1088
$ws.Range("L68").Value = 10000
$ws.Range("M68").Value = -339
$ws.Range("N68").Value = -11498
$ws.Range("H71").Value = 1898.1818
$ws.Range("I71").Value = 1088
$ws.Range("J71").Value = 10000
$ws.Range("K71").Value = 5440
$ws.Range("L71").Value = 50000
$ws.Range("M71").Value = -1696
$ws.Range("N71").Value = -57488
$ws.Range("H80").Value = 38000
$ws.Range("J80").Value = 38000
$ws.Range("L80").Value = 38000
$ws.Range("N80").Value = -40246
$ws.Range("H83").Value = 38000
$ws.Range("J83").Value = 38000
$ws.Range("L83").Value = 114000
$ws.Range("N83").Value = -125232
$ws.Range("H113").Value = 100002360
$ws.Range("I113").Value = 111112620
$ws.Range("K113").Value = 111112620
$ws.Range("M113").Value = -111110450
$ws.Range("H126").Value = 2393.4375
$ws.Range("I126").Value = 1800
$ws.Range("J126").Value = 2855
$ws.Range("K126").Value = 5400
$ws.Range("L126").Value = 8565
$ws.Range("M126").Value = -2930
$ws.Range("N126").Value = -13505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1119.1875
$ws.Range("I113").Value = 454.15384
$ws.Range("J113").Value = 4001
$ws.Range("K113").Value = 1362.46152
$ws.Range("L113").Value = 12003
$ws.Range("M113").Value = 807.5384799999999
$ws.Range("N113").Value = -16343
$ws.Range("H126").Value = 4763891.5
$ws.Range("I126").Value = 1539.7778
$ws.Range("J126").Value = 33338002
$ws.Range("K126").Value = 4619.3334
$ws.Range("L126").Value = 100014006
$ws.Range("M126").Value = -2149.3334
$ws.Range("N126").Value = -100018946
$ws.Range("H132").Value = 10771.714
$ws.Range("I132").Value = 2789.6667
$ws.Range("J132").Value = 30726.834
$ws.Range("K132").Value = 8369.000100000001
$ws.Range("L132").Value = 92180.50199999999
$ws.Range("M132").Value = -5839.000100000001
$ws.Range("N132").Value = -97240.50199999999
$ws.Range("H136").Value = 2217.389
$ws.Range("I136").Value = 1426
$ws.Range("J136").Value = 2850.5
$ws.Range("K136").Value = 4278
$ws.Range("L136").Value = 8551.5
$ws.Range("M136").Value = -1728
$ws.Range("N136").Value = -13651.5

Write-Host "Edits applied successfully"
